$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F
$updates = @{
    3  = 7390
    4  = 7273
    11 = 125
    12 = 210
    13 = 87
    15 = 475
    17 = 22
    19 = 130
    20 = 70
}

# Both "展览" (sheet1) and "全部类型" (sheet4) contain the same table and need
# the same update applied to column F ("想去人数").
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
